$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "fuel_boiler_efficiency" in F1
$ws.Range("F1").Value = "fuel_boiler_efficiency"

# Fill the new column with fuel boiler efficiency values for each boiler row
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0.9
$ws.Range("F4").Value = 0.9
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 0.9

# Update selection to match the last edited cell
$ws.Range("F6").Select()
